$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$c = $ws.Range("A11")
$c.Value = "Pullover"
$c.Font.Family = 1
$c.Font.Size = 10
$c.Font.Name = "Times New Roman"
